# Updated to newest content
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update header text: "C code Programming Ease" -> "Programming Ease" (H2) ---
$ws.Range("H2").Value = "Programming Ease"

# --- Update weights in row 1 ---
$ws.Range("C1").Value = 0.1
$ws.Range("E1").Value = 0.15

# --- Update Arduino row (row 3) ---
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 5

# --- Update Raspberry Pi row (row 4) ---
$ws.Range("C4").Value = 3

# --- Update Tiva row (row 5) ---
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 3
$ws.Range("G5").Value = 4
$ws.Range("H5").Value = 4

# --- Move selection ---
$ws.Range("B15").Select()

# --- Adjust column widths to match target best-fit widths ---
# (New narrow helper columns C, E, F, I appear and A/H narrow because the
# longest entries in those columns got shorter/removed.)
$ws.Columns.Item(1).ColumnWidth = 12.76
$ws.Columns.Item(3).ColumnWidth = 3.76
$ws.Columns.Item(5).ColumnWidth = 4.26
$ws.Columns.Item(6).ColumnWidth = 4.92
$ws.Columns.Item(8).ColumnWidth = 15.26
$ws.Columns.Item(9).ColumnWidth = 8.92
